$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for columns C:F across rows 2-8
# (use Value2 for reads — Value's getter is unreliable on this host)
$orig = @{}
for ($r = 2; $r -le 8; $r++) {
    $orig[$r] = @(
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2
    )
}

# Target layout: rows get the runs/balls/fours/sixes that used to belong
# to a different delivery/row (a reshuffle of the per-innings stat rows).
# Mapping is new row -> source (original) row. Rows 6 and 8 are omitted
# since they map to themselves (no change required).
$map = @{
    2 = 7
    3 = 5
    4 = 3
    5 = 2
    7 = 4
}

foreach ($newRow in $map.Keys) {
    $srcRow = $map[$newRow]
    $vals = $orig[$srcRow]

    $cC = $ws.Cells.Item($newRow, 3)
    $cD = $ws.Cells.Item($newRow, 4)
    $cE = $ws.Cells.Item($newRow, 5)
    $cF = $ws.Cells.Item($newRow, 6)

    # Force text format so the numeric-looking stats are stored as text,
    # matching the workbook's existing "numberStoredAsText" convention.
    $cC.NumberFormat = "@"
    $cD.NumberFormat = "@"
    $cE.NumberFormat = "@"
    $cF.NumberFormat = "@"

    $cC.Value = [string]$vals[0]
    $cD.Value = [string]$vals[1]
    $cE.Value = [string]$vals[2]
    $cF.Value = [string]$vals[3]
}
